# ChartSampleJa.xlsx - valueObject sheet edit
# Commit: i18n and validation issues are updated.
#   * labels should be defined as property style but class getter/setter.
#   * i18n instance is now kept in ValidateConfig class statically.
#
# Concretely this inserts one new row (row 18) into the "valueObject"
# sheet, just below the existing "インタフェイス" (interface) row and
# above the row that used to be row 18 ("ファイナル" / final). The new
# row documents a new "ラベル" (label) flag, flows the "○" marker into
# column C (same as the other boolean-style rows), and records an
# explanatory comment in column D. Everything below row 17 shifts down
# by one row as a natural consequence of the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("valueObject")

# --- Insert the new row -----------------------------------------------
# Inserting at row 18 pushes the previous row 18 ("ファイナル") and
# everything after it down by one row, which is exactly what the diff
# shows (rows 18-90 -> 19-91, formulas/validations/merges renumbered).
$ws.Rows("18").Insert()

# Carry the formatting from row 17 (the row immediately above, which is
# the same "label flag" row template: merged A:B header cell, C column
# list-validated marker, D column free text) down onto the freshly
# inserted row 18.
$ws.Range("A17:D17").Copy()
$ws.Range("A18:D18").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Fill in the new row's content -------------------------------------
$ws.Range("A18").Value = "ラベル"
$ws.Range("C18").Value = "○"
$ws.Range("D18").Value = "/* TypeScript 独自。インタフェイス指定が優先します。 */"

# A18/B18 are merged, matching the other header-style rows (A17:B17,
# A22:B22, ...).
$ws.Range("A18:B18").Merge()

# --- Data validation housekeeping ---------------------------------------
# The "isAbstract" list validation used to cover C16:C17; the newly
# inserted row 18 belongs to the same boolean-flag block, so the range
# grows to C16:C18.
$ws.Range("C16:C17").Validation.Delete()
$ws.Range("C16:C18").Validation.Add(3, 1, 1, "isAbstract", "0")

# --- Fix up the stray _FilterDatabase defined name -----------------------
# This hidden name is left over from a previous AutoFilter and isn't
# touched automatically by the row insert, so shift it by hand (same +1
# row offset as everything else below row 18).
$names = $wb.Names
$fd = $names.Item("valueObject!_FilterDatabase")
$fd.RefersTo = "=valueObject!`$B`$52:`$B`$74"

# --- Restore the view/selection state -----------------------------------
$ws.Activate() | Out-Null
$ws.Range("D19").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 6
